$wb = $excel.ActiveWorkbook

# Add the new "Datasets" header for column A (row 1) on every sheet,
# copying the header formatting already used by the other column headers.
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("A1").Value = "Datasets"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Refresh the recomputed GroupTC-HS figures (column K) on each sheet.

$ws = $wb.Worksheets.Item(1)
$ws.Range("K2").Value = 74.85670399999999
$ws.Range("K3").Value = 82.019752
$ws.Range("K4").Value = 82.04433400000001
$ws.Range("K5").Value = 72.361431
$ws.Range("K6").Value = 91.99444099999999
$ws.Range("K7").Value = 78.522216
$ws.Range("K8").Value = 70.56899199999999
$ws.Range("K9").Value = 81.33661499999999
$ws.Range("K10").Value = 83.42722999999999
$ws.Range("K11").Value = 89.822171
$ws.Range("K12").Value = 93.495626
$ws.Range("K13").Value = 87.142172
$ws.Range("K14").Value = 87.762428
$ws.Range("K15").Value = 72.691294
$ws.Range("K16").Value = 88.75796699999999
$ws.Range("K17").Value = 92.209141
$ws.Range("K18").Value = 88.54284699999999
$ws.Range("K19").Value = 72.695441
$ws.Range("K20").Value = 86.391222
$ws.Range("K21").Value = 93.209444

$ws = $wb.Worksheets.Item(2)
$ws.Range("K2").Value = 1767574
$ws.Range("K3").Value = 363833
$ws.Range("K4").Value = 746650
$ws.Range("K5").Value = 344216
$ws.Range("K6").Value = 8156188
$ws.Range("K7").Value = 5747718
$ws.Range("K8").Value = 10288172
$ws.Range("K9").Value = 11918483
$ws.Range("K10").Value = 5202539
$ws.Range("K11").Value = 19783876
$ws.Range("K12").Value = 149114869
$ws.Range("K13").Value = 53448669
$ws.Range("K14").Value = 84606060
$ws.Range("K16").Value = 1866365486
$ws.Range("K17").Value = 487374950
$ws.Range("K18").Value = 539686341
$ws.Range("K20").Value = 16730350543
$ws.Range("K21").Value = 13066545191

$ws = $wb.Worksheets.Item(3)
$ws.Range("K2").Value = 3.07196
$ws.Range("K3").Value = 9.49371
$ws.Range("K4").Value = 9.963347000000001
$ws.Range("K5").Value = 11.31398
$ws.Range("K6").Value = 4.378134
$ws.Range("K7").Value = 5.107587
$ws.Range("K8").Value = 3.571611
$ws.Range("K9").Value = 5.011336
$ws.Range("K10").Value = 11.666833
$ws.Range("K11").Value = 7.139333
$ws.Range("K12").Value = 3.550772
$ws.Range("K13").Value = 4.69953
$ws.Range("K14").Value = 4.183846
$ws.Range("K15").Value = 12.882196
$ws.Range("K16").Value = 3.116137
$ws.Range("K17").Value = 3.55288
$ws.Range("K18").Value = 3.839044
$ws.Range("K19").Value = 11.34649
$ws.Range("K20").Value = 3.550013
$ws.Range("K21").Value = 2.991866

$ws = $wb.Worksheets.Item(4)
$ws.Range("K2").Value = 50.430386
$ws.Range("K3").Value = 32.782666
$ws.Range("K4").Value = 32.59878
$ws.Range("K5").Value = 28.509947
$ws.Range("K6").Value = 67.370031
$ws.Range("K7").Value = 49.883639
$ws.Range("K8").Value = 42.788118
$ws.Range("K9").Value = 52.005437
$ws.Range("K10").Value = 31.196084
$ws.Range("K11").Value = 53.706665
$ws.Range("K12").Value = 70.215619
$ws.Range("K13").Value = 54.567224
$ws.Range("K14").Value = 56.021738
$ws.Range("K15").Value = 25.469338
$ws.Range("K16").Value = 74.01255
$ws.Range("K17").Value = 66.44068900000001
$ws.Range("K18").Value = 65.017385
$ws.Range("K19").Value = 28.798812
$ws.Range("K20").Value = 77.73144499999999
$ws.Range("K21").Value = 68.312934
